# Append 10 new DQS log rows (71-80) to Sheet1, matching the new runs
# logged on 2024-12-03 for CK_CU_BOUNDARY_En and Conservation_Unit_Data_20220902.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colAllCols = "ACT_ID, ANALYSIS_YR, STREAM_ID, SPL_ID, NATURAL_ADULT_SPAWNERS, NATURAL_JACK_SPAWNERS, NATURAL_SPAWNERS_TOTAL, ADULT_BROODSTOCK_REMOVALS, JACK_BROODSTOCK_REMOVALS, TOTAL_BROODSTOCK_REMOVALS, OTHER_REMOVALS, TOTAL_RETURN_TO_RIVER, UNSPECIFIED_RETURN, NO_INSPECTIONS_USED, MAX_ESTIMATE, EFFECTIVE_FEMALES, WEIGHTED_PCT_SPAWN, OTHER_ADULT_REMOVALS, OTHER_JACK_REMOVALS, TOT_ADULT_RET_RIVER, TOT_JACK_RET_RIVER, JUV_PRES_TYP, POP_ID, SBJ_ID"

# Row 71
$ws.Cells.Item(71, 1).Value = "CK_CU_BOUNDARY_En"
$ws.Cells.Item(71, 3).Value = "Accuracy (A1)"
$ws.Cells.Item(71, 4).Value = "2024-12-03 12:06:55"
$ws.Cells.Item(71, 5).Value = "no threshold"
$ws.Cells.Item(71, 7).Value = "OnakD"

# Row 72
$ws.Cells.Item(72, 1).Value = "CK_CU_BOUNDARY_En"
$ws.Cells.Item(72, 2).Value = "All columns"
$ws.Cells.Item(72, 3).Value = "Accuracy (A3)"
$ws.Cells.Item(72, 4).Value = "2024-12-03 12:06:56"
$ws.Cells.Item(72, 5).Value = "no threshold"
$ws.Cells.Item(72, 6).Value = 1
$ws.Cells.Item(72, 7).Value = "OnakD"

# Row 73
$ws.Cells.Item(73, 1).Value = "CK_CU_BOUNDARY_En"
$ws.Cells.Item(73, 2).Value = "All columns"
$ws.Cells.Item(73, 3).Value = "Completeness (P)"
$ws.Cells.Item(73, 4).Value = "2024-12-03 12:06:56"
$ws.Cells.Item(73, 5).Value = 0.75
$ws.Cells.Item(73, 6).Value = 1
$ws.Cells.Item(73, 7).Value = "OnakD"

# Row 74
$ws.Cells.Item(74, 1).Value = "CK_CU_BOUNDARY_En"
$ws.Cells.Item(74, 3).Value = "Accuracy (A1)"
$ws.Cells.Item(74, 4).Value = "2024-12-03 12:13:31"
$ws.Cells.Item(74, 5).Value = "no threshold"
$ws.Cells.Item(74, 7).Value = "OnakD"

# Row 75
$ws.Cells.Item(75, 1).Value = "CK_CU_BOUNDARY_En"
$ws.Cells.Item(75, 2).Value = "All columns"
$ws.Cells.Item(75, 3).Value = "Accuracy (A3)"
$ws.Cells.Item(75, 4).Value = "2024-12-03 12:13:31"
$ws.Cells.Item(75, 5).Value = "no threshold"
$ws.Cells.Item(75, 6).Value = 1
$ws.Cells.Item(75, 7).Value = "OnakD"

# Row 76
$ws.Cells.Item(76, 1).Value = "CK_CU_BOUNDARY_En"
$ws.Cells.Item(76, 2).Value = "All columns"
$ws.Cells.Item(76, 3).Value = "Completeness (P)"
$ws.Cells.Item(76, 4).Value = "2024-12-03 12:13:31"
$ws.Cells.Item(76, 5).Value = 0.75
$ws.Cells.Item(76, 6).Value = 1
$ws.Cells.Item(76, 7).Value = "OnakD"

# Row 77
$ws.Cells.Item(77, 1).Value = "Conservation_Unit_Data_20220902"
$ws.Cells.Item(77, 2).Value = $colAllCols
$ws.Cells.Item(77, 3).Value = "Accuracy (A1)"
$ws.Cells.Item(77, 4).Value = "2024-12-03 12:15:28"
$ws.Cells.Item(77, 5).Value = "no threshold"
$ws.Cells.Item(77, 6).Value = 0.9968144750254843
$ws.Cells.Item(77, 7).Value = "OnakD"

# Row 78
$ws.Cells.Item(78, 1).Value = "Conservation_Unit_Data_20220902"
$ws.Cells.Item(78, 2).Value = $colAllCols
$ws.Cells.Item(78, 3).Value = "Accuracy (A2)"
$ws.Cells.Item(78, 4).Value = "2024-12-03 12:15:29"
$ws.Cells.Item(78, 5).Value = 0.85
$ws.Cells.Item(78, 6).Value = 1
$ws.Cells.Item(78, 7).Value = "OnakD"

# Row 79
$ws.Cells.Item(79, 1).Value = "Conservation_Unit_Data_20220902"
$ws.Cells.Item(79, 2).Value = "All columns"
$ws.Cells.Item(79, 3).Value = "Accuracy (A3)"
$ws.Cells.Item(79, 4).Value = "2024-12-03 12:15:31"
$ws.Cells.Item(79, 5).Value = "no threshold"
$ws.Cells.Item(79, 6).Value = 1
$ws.Cells.Item(79, 7).Value = "OnakD"

# Row 80
$ws.Cells.Item(80, 1).Value = "Conservation_Unit_Data_20220902"
$ws.Cells.Item(80, 2).Value = "All columns"
$ws.Cells.Item(80, 3).Value = "Completeness (P)"
$ws.Cells.Item(80, 4).Value = "2024-12-03 12:15:32"
$ws.Cells.Item(80, 5).Value = 0.75
$ws.Cells.Item(80, 6).Value = 0.9240346358763629
$ws.Cells.Item(80, 7).Value = "OnakD"
